$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.019.77"
$ws.Range("E2").Value = "  -0.98%  "

# Row 3
$ws.Range("D3").Value = "1.619.85"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "

# Row 6
$ws.Range("E6").Value = "  -0.06%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.81%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

# Row 12
$ws.Range("D12").Value = "1.622.93"
$ws.Range("E12").Value = "  -1.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.539"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.72%  "

# Row 16
$ws.Range("D16").Value = "27.004.13"
$ws.Range("E16").Value = "  -1.04%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0748"
$ws.Range("E17").Value = "  +0.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.72%  "

# Row 19
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "

# Row 22
$ws.Range("E22").Value = "  -5.81%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "

# Row 27
$ws.Range("E27").Value = "  -1.99%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.21%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0516"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.82%  "

# Row 30
$ws.Range("E30").Value = "  -1.28%  "

# Row 31
$ws.Range("E31").Value = "  -0.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.745"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +34.89%  "

# Row 33
$ws.Range("E33").Value = "  -0.25%  "

# Row 34
$ws.Range("D34").Value = "1.345.29"
$ws.Range("E34").Value = "  +3.23%  "

# Row 35
$ws.Range("E35").Value = "  -0.89%  "

# Row 36
$ws.Range("E36").Value = "  -0.44%  "

# Row 37
$ws.Range("E37").Value = "  +1.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.851"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.18%  "

# Row 39
$ws.Range("E39").Value = "  -0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.800"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.30%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").Value = "1.756.92"
$ws.Range("E44").Value = "  -1.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.50%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.870"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +29.64%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.14%  "

# Row 48
$ws.Range("E48").Value = "  -0.01%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.51%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.84%  "
